# Update addresses for three pubs and re-apply the default font across the
# data rows (mirrors the formatting refresh + address lookups added upstream),
# then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the three missing addresses (new shared strings) ---------------
$ws.Range("B12").Value = "Vaníčkova Panská, 400 01 Ústí nad Labem-město"
$ws.Range("B13").Value = "Hrnčířská 10, 400 01 Ústí nad Labem-město"
$ws.Range("B14").Value = "Pivovarská 3380/5, 400 01 Ústí nad Labem-centrum"

# Give the three new address cells their own font (Helvetica Neue, dark grey)
$ws.Range("B12:B14").Font.Color = 2367776
$ws.Range("B12:B14").Font.Name = "Helvetica Neue"

# --- Re-apply (stamp) the base font on the rest of the data cells -----------
# (covers A2:E7, C8:D8, A9:E11, A12:A14, C12:E12, C13:E13, C14:E14 -
#  i.e. every data cell except the header row and the A8/B8/E8 cells that
#  already carry their own explicit style)
$ws.Range("A2:E7").Font.Name = "Aptos Narrow"
$ws.Range("C8:D8").Font.Name = "Aptos Narrow"
$ws.Range("A9:E11").Font.Name = "Aptos Narrow"
$ws.Range("A12:A14").Font.Name = "Aptos Narrow"
$ws.Range("C12:E12").Font.Name = "Aptos Narrow"
$ws.Range("C13:E13").Font.Name = "Aptos Narrow"
$ws.Range("C14:E14").Font.Name = "Aptos Narrow"

# --- Move the active selection ----------------------------------------------
$ws.Range("D21").Select()
